$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1299.75
$ws.Range("I18").Value = 1100
$ws.Range("K18").Value = 1100
$ws.Range("M18").Value = -816
$ws.Range("H32").Value = 5914.846
$ws.Range("J32").Value = 5302.3335
$ws.Range("L32").Value = 5302.3335
$ws.Range("N32").Value = -5954.3335
$ws.Range("H40").Value = 4263
$ws.Range("I40").Value = 3482.2942
$ws.Range("K40").Value = 3482.2942
$ws.Range("M40").Value = -3307.2942
$ws.Range("H55").Value = 512.4666999999999
$ws.Range("J55").Value = 1027.1666
$ws.Range("L55").Value = 1027.1666
$ws.Range("N55").Value = -1455.1666
$ws.Range("H86").Value = 3707.5757
$ws.Range("I86").Value = 3034.8125
$ws.Range("J86").Value = 4340.7646
$ws.Range("K86").Value = 3034.8125
$ws.Range("L86").Value = 4340.7646
$ws.Range("M86").Value = -1911.8125
$ws.Range("N86").Value = -6586.7646
$ws.Range("H89").Value = 3707.5757
$ws.Range("I89").Value = 3034.8125
$ws.Range("J89").Value = 4340.7646
$ws.Range("K89").Value = 15174.0625
$ws.Range("L89").Value = 21703.823
$ws.Range("M89").Value = -9558.0625
$ws.Range("N89").Value = -32935.823
$ws.Range("H98").Value = 261794.27
$ws.Range("I98").Value = 1039.8422
$ws.Range("K98").Value = 1039.8422
$ws.Range("M98").Value = 458.1578
$ws.Range("H122").Value = 261794.27
$ws.Range("I122").Value = 1039.8422
$ws.Range("K122").Value = 3119.5266
$ws.Range("M122").Value = -669.5266000000001
$ws.Range("H129").Value = 2140
$ws.Range("I129").Value = 950
$ws.Range("J129").Value = 2933.3333
$ws.Range("K129").Value = 2850
$ws.Range("L129").Value = 8799.999899999999
$ws.Range("M129").Value = 2150
$ws.Range("N129").Value = -18799.9999
$ws.Range("H138").Value = 2773.281
$ws.Range("I138").Value = 1813.6666
$ws.Range("K138").Value = 5440.9998
$ws.Range("M138").Value = -300.9997999999996
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3825.7932
$ws.Range("I2").Value = 314.20834
$ws.Range("K2").Value = 314.20834
$ws.Range("M2").Value = -201.20834
$ws.Range("H74").Value = 16670945
$ws.Range("I74").Value = 25643606
$ws.Range("K74").Value = 25643606
$ws.Range("M74").Value = -25642732
$ws.Range("H77").Value = 16670945
$ws.Range("I77").Value = 25643606
$ws.Range("K77").Value = 128218030
$ws.Range("M77").Value = -128213662
$ws.Range("H116").Value = 3825.7932
$ws.Range("I116").Value = 314.20834
$ws.Range("K116").Value = 314.20834
$ws.Range("M116").Value = 1979.79166
$ws.Range("H132").Value = 2381.0303
$ws.Range("I132").Value = 1806.7587
$ws.Range("J132").Value = 6544.5
$ws.Range("K132").Value = 5420.2761
$ws.Range("L132").Value = 19633.5
$ws.Range("M132").Value = -2890.2761
$ws.Range("N132").Value = -24693.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3825.7932
$ws.Range("I3").Value = 314.20834
$ws.Range("K3").Value = 314.20834
$ws.Range("M3").Value = -200.20834
$ws.Range("H22").Value = 519.5
$ws.Range("I22").Value = 522.2857
$ws.Range("K22").Value = 522.2857
$ws.Range("M22").Value = -349.2857
$ws.Range("H99").Value = 1932.8182
$ws.Range("I99").Value = 1508.8572
$ws.Range("J99").Value = 2674.75
$ws.Range("K99").Value = 1508.8572
$ws.Range("L99").Value = 2674.75
$ws.Range("M99").Value = -10.85719999999992
$ws.Range("N99").Value = -5670.75
$ws.Range("H134").Value = 2407.8572
$ws.Range("I134").Value = 2515.077
$ws.Range("K134").Value = 7545.231000000001
$ws.Range("M134").Value = -5010.231000000001
$ws.Range("H138").Value = 65030.285
$ws.Range("J138").Value = 65030.285
$ws.Range("L138").Value = 65030.285
$ws.Range("N138").Value = -75310.285
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2149.4443
$ws.Range("I22").Value = 608
$ws.Range("K22").Value = 608
$ws.Range("M22").Value = -258
$ws.Range("H99").Value = 3533.2727
$ws.Range("I99").Value = 3149.25
$ws.Range("K99").Value = 3149.25
$ws.Range("M99").Value = -1651.25
$ws.Range("H105").Value = 3237.6155
$ws.Range("I105").Value = 1384.875
$ws.Range("J105").Value = 6202
$ws.Range("K105").Value = 1384.875
$ws.Range("L105").Value = 6202
$ws.Range("M105").Value = 362.125
$ws.Range("N105").Value = -9696
$ws.Range("H126").Value = 3533.2727
$ws.Range("I126").Value = 3149.25
$ws.Range("K126").Value = 9447.75
$ws.Range("M126").Value = -6977.75
$ws.Range("H132").Value = 3567.5881
$ws.Range("I132").Value = 2638.739
$ws.Range("K132").Value = 7916.217000000001
$ws.Range("M132").Value = -5386.217000000001
$ws.Range("H134").Value = 2119.162
$ws.Range("I134").Value = 1366.1923
$ws.Range("K134").Value = 4098.5769
$ws.Range("M134").Value = -1563.5769
$ws.Range("H135").Value = 69178
$ws.Range("J135").Value = 69178
$ws.Range("L135").Value = 69178
$ws.Range("N135").Value = -79318
$ws.Range("H138").Value = 69284
$ws.Range("J138").Value = 69284
$ws.Range("L138").Value = 69284
$ws.Range("N138").Value = -79564
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 274376.12
$ws.Range("J141").Value = 274376.12
$ws.Range("L141").Value = 274376.12
$ws.Range("N141").Value = -284736.12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 146.09091
$ws.Range("J23").Value = 212.42857
$ws.Range("L23").Value = 637.28571
$ws.Range("N23").Value = -1107.28571
$ws.Range("H136").Value = 3975.8
$ws.Range("I136").Value = 3719.75
$ws.Range("K136").Value = 11159.25
$ws.Range("M136").Value = -6059.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H80").Value = 504571.6
$ws.Range("J80").Value = 6361.2
$ws.Range("L80").Value = 6361.2
$ws.Range("N80").Value = -8357.200000000001
$ws.Range("H83").Value = 504571.6
$ws.Range("J83").Value = 6361.2
$ws.Range("L83").Value = 31806
$ws.Range("N83").Value = -41790
$ws.Range("H102").Value = 2192.6724
$ws.Range("I102").Value = 1475.3112
$ws.Range("K102").Value = 1475.3112
$ws.Range("M102").Value = 146.6887999999999
$ws.Range("H113").Value = 2693.1155
$ws.Range("I113").Value = 2223.625
$ws.Range("K113").Value = 2223.625
$ws.Range("M113").Value = -53.625
$ws.Range("H122").Value = 10464.909
$ws.Range("I122").Value = 14146.091
$ws.Range("K122").Value = 42438.273
$ws.Range("M122").Value = -39988.273
$ws.Range("H132").Value = 2113.7666
$ws.Range("I132").Value = 1403.4286
$ws.Range("K132").Value = 4210.2858
$ws.Range("M132").Value = -1680.2858
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2941.5
$ws.Range("I136").Value = 1458.3684
$ws.Range("K136").Value = 4375.1052
$ws.Range("M136").Value = -1825.1052
